$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 463.83334
$ws.Range("I9").Value = 356.8
$ws.Range("J9").Value = 999
$ws.Range("K9").Value = 356.8
$ws.Range("L9").Value = 999
$ws.Range("M9").Value = -187.8
$ws.Range("N9").Value = -1337
$ws.Range("H17").Value = 1797.5
$ws.Range("J17").Value = 1797.5
$ws.Range("L17").Value = 5392.5
$ws.Range("N17").Value = -5728.5
$ws.Range("H32").Value = 9499.5
$ws.Range("I32").Value = 7999
$ws.Range("K32").Value = 7999
$ws.Range("M32").Value = -7673
$ws.Range("H40").Value = 6659.84
$ws.Range("I40").Value = 6998
$ws.Range("J40").Value = 6645.75
$ws.Range("K40").Value = 6998
$ws.Range("L40").Value = 6645.75
$ws.Range("M40").Value = -6823
$ws.Range("N40").Value = -6995.75
$ws.Range("H69").Value = 6012
$ws.Range("J69").Value = 5761.25
$ws.Range("L69").Value = 17283.75
$ws.Range("N69").Value = -19031.75
$ws.Range("H72").Value = 6012
$ws.Range("J72").Value = 5761.25
$ws.Range("L72").Value = 51851.25
$ws.Range("N72").Value = -60587.25
$ws.Range("H127").Value = 887.5
$ws.Range("I127").Value = 1187.5
$ws.Range("J127").Value = 587.5
$ws.Range("K127").Value = 3562.5
$ws.Range("L127").Value = 1762.5
$ws.Range("M127").Value = 1397.5
$ws.Range("N127").Value = -11682.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3625
$ws.Range("I61").Value = 3625
$ws.Range("K61").Value = 3625
$ws.Range("M61").Value = -3413
$ws.Range("H97").Value = 2479.3333
$ws.Range("I97").Value = 1431.5
$ws.Range("K97").Value = 1431.5
$ws.Range("M97").Value = -935.5
$ws.Range("H122").Value = 2253
$ws.Range("I122").Value = 2029.5
$ws.Range("K122").Value = 6088.5
$ws.Range("M122").Value = -3638.5
$ws.Range("H132").Value = 2820.125
$ws.Range("I132").Value = 2426.8333
$ws.Range("K132").Value = 7280.499899999999
$ws.Range("M132").Value = -4750.499899999999
$ws.Range("H136").Value = 3625
$ws.Range("I136").Value = 3625
$ws.Range("K136").Value = 10875
$ws.Range("M136").Value = -8325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""
$ws.Range("H134").Value = 4836.222
$ws.Range("I134").Value = 3004.3333
$ws.Range("K134").Value = 9012.999899999999
$ws.Range("M134").Value = -6477.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19713
$ws.Range("H31").Value = 3997.7646
$ws.Range("I31").Value = 1913.4286
$ws.Range("K31").Value = 1913.4286
$ws.Range("M31").Value = -1618.4286
$ws.Range("H34").Value = 3997.7646
$ws.Range("I34").Value = 1913.4286
$ws.Range("K34").Value = 1913.4286
$ws.Range("M34").Value = -1711.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 20005360
$ws.Range("I11").Value = 20005360
$ws.Range("K11").Value = 60016080
$ws.Range("M11").Value = -60015940
$ws.Range("H39").Value = 3525.5
$ws.Range("I39").Value = 2548.3333
$ws.Range("J39").Value = 4502.6665
$ws.Range("K39").Value = 7644.999899999999
$ws.Range("L39").Value = 13507.9995
$ws.Range("M39").Value = -7350.999899999999
$ws.Range("N39").Value = -14095.9995
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 1899.5
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 5698.5
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -7320.5
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 1899.5
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 17095.5
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -25207.5
$ws.Range("H110").Value = 3500
$ws.Range("J110").Value = 3500
$ws.Range("L110").Value = 10500
$ws.Range("N110").Value = -18680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2406.5
$ws.Range("I102").Value = 2406.5
$ws.Range("K102").Value = 2406.5
$ws.Range("M102").Value = -784.5
$ws.Range("H132").Value = 8221.556
$ws.Range("I132").Value = 7749.1665
$ws.Range("K132").Value = 23247.4995
$ws.Range("M132").Value = -20717.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8249.5
$ws.Range("I7").Value = 8332.666999999999
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 8332.666999999999
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -8220.666999999999
$ws.Range("N7").Value = -8224
$ws.Range("H40").Value = 9000
$ws.Range("I40").Value = 9000
$ws.Range("K40").Value = 9000
$ws.Range("M40").Value = -8864
$ws.Range("H100").Value = 7870.3335
$ws.Range("I100").Value = 10555.5
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 10555.5
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -10014.5
$ws.Range("N100").Value = -3582
$ws.Range("H126").Value = 8249.5
$ws.Range("I126").Value = 8332.666999999999
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 24998.001
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -22528.001
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 30380.6
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 60761.2
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -59700.2
$ws.Range("N81").Value = -4122
$ws.Range("I84").Value = 30380.6
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 303806
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -298502
$ws.Range("N84").Value = -20608
$ws.Range("H100").Value = 748.125
$ws.Range("J100").Value = 485
$ws.Range("L100").Value = 970
$ws.Range("N100").Value = -2052
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530
